# Auto-generated edit script applying cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.363.34'
$ws.Range("E2").Value = '  -0.30%  '

$ws.Range("D3").Value = '1.847.56'
$ws.Range("E3").Value = '  -0.29%  '

$ws.Range("D4").Value = '''0.9995'
$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").Value = '''240.24'
$ws.Range("E5").Value = '  -0.25%  '

$ws.Range("D6").Value = '''0.6281'
$ws.Range("E6").Value = '  -0.46%  '

$ws.Range("D8").Value = '''0.07601'
$ws.Range("E8").Value = '  -0.92%  '

$ws.Range("D9").Value = '''0.2912'
$ws.Range("E9").Value = '  -1.12%  '

$ws.Range("D10").Value = '''24.55'
$ws.Range("E10").Value = '  -0.56%  '

$ws.Range("D11").Value = '''0.07746'
$ws.Range("E11").Value = '  -0.10%  '

$ws.Range("D12").Value = '1.847.91'
$ws.Range("E12").Value = '  -0.26%  '

$ws.Range("E13").Value = '  -0.34%  '

$ws.Range("D14").Value = '''0.6794'
$ws.Range("E14").Value = '  -0.34%  '

$ws.Range("D15").Value = '''0.00001042'
$ws.Range("E15").Value = '  -4.12%  '

$ws.Range("D16").Value = '''82.97'

$ws.Range("D17").Value = '''6.113'
$ws.Range("E17").Value = '  -0.95%  '

$ws.Range("D18").Value = '29.368.09'
$ws.Range("E18").Value = '  -0.35%  '

$ws.Range("D19").Value = '''228.48'
$ws.Range("E19").Value = '  -0.51%  '

$ws.Range("E20").Value = '  -1.27%  '

$ws.Range("D21").Value = '''1.000'
$ws.Range("E21").Value = '  -0.04%  '

$ws.Range("D22").Value = '''7.434'
$ws.Range("E22").Value = '  -0.32%  '

$ws.Range("E23").Value = '  -0.03%  '

$ws.Range("D24").Value = '''158.86'
$ws.Range("E24").Value = '  +1.25%  '

$ws.Range("D25").Value = '''0.1390'
$ws.Range("E25").Value = '  +0.30%  '

$ws.Range("D26").Value = '''8.435'
$ws.Range("E26").Value = '  +0.23%  '

$ws.Range("D27").Value = '''17.65'
$ws.Range("E27").Value = '  -0.43%  '

$ws.Range("D28").Value = '''1.432'
$ws.Range("E28").Value = '  +8.02%  '

$ws.Range("E29").Value = '  -0.19%  '

$ws.Range("D30").Value = '''0.05623'
$ws.Range("E30").Value = '  -1.17%  '

$ws.Range("D31").Value = '''4.107'

$ws.Range("D32").Value = '''4.035'
$ws.Range("E32").Value = '  -0.44%  '

$ws.Range("D33").Value = '''1.825'
$ws.Range("E33").Value = '  -1.45%  '

$ws.Range("D34").Value = '''1.156'
$ws.Range("E34").Value = '  -0.66%  '

$ws.Range("D35").Value = '''0.6960'
$ws.Range("E35").Value = '  -1.33%  '

$ws.Range("D36").Value = '''2.581'
$ws.Range("E36").Value = '  -0.16%  '

$ws.Range("D37").Value = '''0.01830'
$ws.Range("E37").Value = '  +2.00%  '

$ws.Range("D38").Value = '1.234.50'
$ws.Range("E38").Value = '  +1.30%  '

$ws.Range("D39").Value = '''2.728'
$ws.Range("E39").Value = '  -2.01%  '

$ws.Range("D40").Value = '''6.397'
$ws.Range("E40").Value = '  -2.30%  '

$ws.Range("D41").Value = '''0.8988'
$ws.Range("E41").Value = '  -1.25%  '

$ws.Range("D42").Value = '''0.9998'
$ws.Range("E42").Value = '  -0.16%  '

$ws.Range("D43").Value = '''101.30'
$ws.Range("E43").Value = '  -0.42%  '

$ws.Range("E44").Value = '  -1.59%  '

$ws.Range("D45").Value = '''7.130'
$ws.Range("E45").Value = '  +0.07%  '

$ws.Range("D46").Value = '''0.3993'
$ws.Range("E46").Value = '  -0.69%  '

$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").Value = '''0.00000000115'
$ws.Range("E47").Value = '  -4.34%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '''8.977'
$ws.Range("E48").Value = '  -0.73%  '

$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").Value = '''0.1149'
$ws.Range("E49").Value = '  +1.20%  '

$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").Value = '''1.678'
$ws.Range("E50").Value = '  -0.49%  '

$ws.Range("D51").Value = '''0.05699'
$ws.Range("E51").Value = '  -0.31%  '

Write-Host "Applied 99 cell updates to cryptos worksheet"
